$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.931.65"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "1.874.56"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("D4").Value = "0.9986"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "306.55"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").Value = "0.9987"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "0.5165"
$ws.Range("E7").Value = "  +1.53%  "
$ws.Range("D8").Value = "0.3707"
$ws.Range("E8").Value = "  +1.49%  "
$ws.Range("D9").Value = "0.07178"
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("D10").Value = "0.8977"
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("D11").Value = "20.71"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07556"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.880.66"
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("D14").Value = "94.69"
$ws.Range("E14").Value = "  +3.97%  "
$ws.Range("D15").Value = "5.242"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").Value = "0.9986"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "0.000008471"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").Value = "14.23"
$ws.Range("E18").Value = "  +1.35%  "
$ws.Range("D19").Value = "0.9983"
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "26.975.02"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").Value = "5.025"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").Value = "2.115.58"
$ws.Range("E22").Value = "  +1.23%  "
$ws.Range("D23").Value = "10.39"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.430"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "145.79"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").Value = "1.782"
$ws.Range("E26").Value = "  -2.11%  "
$ws.Range("D27").Value = "17.99"
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("D28").Value = "2.109"
$ws.Range("E28").Value = "  +3.03%  "
$ws.Range("D29").Value = "114.36"
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("D30").Value = "4.917"
$ws.Range("E30").Value = "  +5.32%  "
$ws.Range("D31").Value = "4.749"
$ws.Range("E31").Value = "  +2.62%  "
$ws.Range("D32").Value = "0.09179"
$ws.Range("E32").Value = "  -0.86%  "
$ws.Range("D33").Value = "0.05032"
$ws.Range("E33").Value = "  -1.42%  "
$ws.Range("D34").Value = "0.7527"
$ws.Range("E34").Value = "  +2.99%  "
$ws.Range("D35").Value = "2.997"
$ws.Range("D36").Value = "1.174"
$ws.Range("E36").Value = "  +2.15%  "
$ws.Range("D37").Value = "3.256"
$ws.Range("E37").Value = "  +2.15%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "2.499"
$ws.Range("E38").Value = "  +1.46%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01989"
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5570"
$ws.Range("E40").Value = "  +5.49%  "
$ws.Range("D42").Value = "6.567"
$ws.Range("E42").Value = "  +1.63%  "
$ws.Range("D43").Value = "116.47"
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("D44").Value = "8.754"
$ws.Range("E44").Value = "  +4.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1500"
$ws.Range("E45").Value = "  +1.89%  "
$ws.Range("D46").Value = "0.4774"
$ws.Range("E46").Value = "  +2.79%  "
$ws.Range("D47").Value = "0.9981"
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("D48").Value = "10.09"
$ws.Range("E48").Value = "  +1.70%  "
$ws.Range("D49").Value = "1.565"
$ws.Range("E49").Value = "  +0.71%  "
$ws.Range("D50").Value = "37.08"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").Value = "63.28"
$ws.Range("E51").Value = "  +0.55%  "
